$d = $word.ActiveDocument

# Identify the paragraphs to remove:
#  - paragraphs that consist solely of an inline picture (the 3 screenshot
#    images that were dropped from the document), and
#  - the empty "separator" paragraphs (no runs, just <w:spacing w:before="40"/>
#    i.e. 2pt space-before) that used to sit right after each code table.
$count = $d.Paragraphs.Count
$toDelete = New-Object System.Collections.ArrayList

for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $ils = $p.Range.InlineShapes.Count
    $sb = $p.Format.SpaceBefore
    $txt = $p.Range.Text
    $isEmptyText = ($txt.Trim("`r") -eq "")

    if ($ils -gt 0) {
        [void]$toDelete.Add($i)
    } elseif ($isEmptyText -and $sb -eq 2) {
        [void]$toDelete.Add($i)
    }
}

# Delete from the bottom up so earlier indices stay valid.
for ($j = $toDelete.Count - 1; $j -ge 0; $j--) {
    $idx = $toDelete[$j]
    $p = $d.Paragraphs.Item($idx)
    $p.Range.Delete()
}

Write-Output ("Deleted paragraphs: " + ($toDelete -join ","))
